$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "sector" column header
$ws.Range("F1").Value = "sector"

# Add sector values for existing rows
$ws.Range("F2").Value = "health"
$ws.Range("F3").Value = "finance"
$ws.Range("F4").Value = "education"

# Add new row 5 - duplicate of row 4 but with sector "energy"
$ws.Range("A5").Value = "abcd2222"
$ws.Range("B5").Value = "abcd2"
$ws.Range("C5").Value = "bse"
$ws.Range("D5").Value = 56
$ws.Range("E5").Value = "2020-12-18T09:00:00.0000000"
$ws.Range("F5").Value = "energy"

# Update companyName and companyCode values for existing rows (2-4)
$ws.Range("A2:A4").Value = "abcd2222"
$ws.Range("B2:B4").Value = "abcd2"

# Update the selection to E3 (to match the diff)
$ws.Range("E3").Select()
